{"js": "// Update the sentencing-entry dates:\n//   \"June 08, 2022\"  -> \"June 09, 2022\"   (appears 3 times)\n//   \"August 07, 2022\" -> \"August 08, 2022\" (appears 1 time)\nconst body = context.document.body;\n\nconst replacements = [\n  [\"June 08, 2022\", \"June 09, 2022\"],\n  [\"August 07, 2022\", \"August 08, 2022\"]\n];\n\nfor (const [findText, replaceText] of replacements) {\n  const results = body.search(findText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(replaceText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update the sentencing-entry dates:\n#   \"June 08, 2022\"   -> \"June 09, 2022\"    (appears 3 times)\n#   \"August 07, 2022\" -> \"August 08, 2022\"  (appears 1 time)\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"June 08, 2022\"\n$find.Replacement.Text = \"June 09, 2022\"\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n\n$find2 = $d.Content.Find\n$find2.Text = \"August 07, 2022\"\n$find2.Replacement.Text = \"August 08, 2022\"\n$find2.Execute($find2.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find2.Replacement.Text, 2)\n"}
